$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.883.69"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").Value = "2.672.11"
$ws.Range("E3").Value = "  +7.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "113.45"
$ws.Range("E5").Value = "  +8.84%  "
$ws.Range("D6").Value = "325.87"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "0.527"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").Value = "40.68"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").Value = "20.09"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "0.0821"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  +4.74%  "
$ws.Range("D15").Value = "3.092.89"
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("D16").Value = "2.679.40"
$ws.Range("E16").Value = "  +7.85%  "
$ws.Range("D17").Value = "0.872"
$ws.Range("E17").Value = "  +5.98%  "
$ws.Range("D18").Value = "49.855.44"
$ws.Range("E18").Value = "  +4.33%  "
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("E22").Value = "  +3.74%  "
$ws.Range("D23").Value = "71.79"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "276.14"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  +4.98%  "
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").Value = "  +6.99%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "36.10"
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("D31").Value = "0.139"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "50.24"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("D34").Value = "19.50"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").Value = "0.0805"
$ws.Range("E35").Value = "  +5.08%  "
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").Value = "  +12.39%  "
$ws.Range("E38").Value = "  +7.26%  "
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +9.66%  "
$ws.Range("D40").Value = "125.71"
$ws.Range("E40").Value = "  +5.01%  "
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").Value = "22.31"
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  +6.72%  "
$ws.Range("D45").Value = "2.122.65"
$ws.Range("E45").Value = "  +6.93%  "
$ws.Range("D46").Value = "3.31"
$ws.Range("E46").Value = "  +7.15%  "
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +10.19%  "
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +7.03%  "
